# Updated the affixes pages, fixed a couple bugs.
#
# Adds a new "affix_type" column (AN) to the Affixes sheet:
#   - AN1 gets the new header "affix_type"
#   - AN2:AN25 all get the value 17 (every existing affix row)
# and updates the sheet's view/selection to point at the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (row 1)
$ws.Range("AN1").Value = "affix_type"

# New column data (rows 2-25)
for ($row = 2; $row -le 25; $row++) {
    $ws.Cells.Item($row, 40).Value = 17
}

# Scroll/select the new column so it becomes the active selection, mirroring
# the author having just finished filling it in.
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 29
$ws.Range("AN2:AN25").Select() | Out-Null
